# Apply the May 9 2024 cryptos data refresh.
#
# The "Price" column (D) stores values as plain text (e.g. "596.75"),
# not numbers. Excel will silently coerce plain-decimal-looking text
# into a numeric value when assigned through .Value, so for any such
# cell we briefly force Text formatting while writing the value and
# then restore the cell style to "Normal" so no visible/style change
# is left behind - only the text content changes, exactly like the
# source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '61.657.52'
$ws.Range("E2").Value = '  -1.49%  '

$ws.Range("D3").Value = '3.010.80'
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("E4").Value = '  -0.09%  '

Set-TextValue "D5" '596.75'
$ws.Range("E5").Value = '  +1.87%  '

Set-TextValue "D6" '144.13'
$ws.Range("E6").Value = '  -2.09%  '

$ws.Range("E7").Value = '  +0.01%  '

Set-TextValue "D8" '0.523'
$ws.Range("E8").Value = '  +0.18%  '

$ws.Range("D9").Value = '3.009.07'
$ws.Range("E9").Value = '  -0.29%  '

$ws.Range("E10").Value = '  -1.33%  '

Set-TextValue "D11" '5.91'
$ws.Range("E11").Value = '  +2.34%  '

Set-TextValue "D12" '0.462'
$ws.Range("E12").Value = '  +4.43%  '

Set-TextValue "D13" '0.0000230'
$ws.Range("E13").Value = '  +0.00%  '

Set-TextValue "D14" '34.43'
$ws.Range("E14").Value = '  -1.23%  '

$ws.Range("E15").Value = '  +2.24%  '

$ws.Range("D16").Value = '3.494.44'
$ws.Range("E16").Value = '  -0.61%  '

Set-TextValue "D17" '7.05'
$ws.Range("E17").Value = '  +0.26%  '

$ws.Range("D18").Value = '61.613.99'
$ws.Range("E18").Value = '  -1.46%  '

$ws.Range("D19").Value = '3.002.59'
$ws.Range("E19").Value = '  -0.47%  '

Set-TextValue "D20" '454.83'
$ws.Range("E20").Value = '  -2.07%  '

Set-TextValue "D21" '14.05'
$ws.Range("E21").Value = '  +0.69%  '

Set-TextValue "D22" '0.689'
$ws.Range("E22").Value = '  +0.33%  '

Set-TextValue "D23" '7.38'
$ws.Range("E23").Value = '  +0.24%  '

Set-TextValue "D24" '82.09'
$ws.Range("E24").Value = '  +2.66%  '

Set-TextValue "D25" '2.24'
$ws.Range("E25").Value = '  -4.15%  '

Set-TextValue "D26" '10.66'
$ws.Range("E26").Value = '  +4.83%  '

Set-TextValue "D27" '12.00'
$ws.Range("E27").Value = '  -3.26%  '

$ws.Range("E28").Value = '  +0.17%  '

Set-TextValue "D29" '2.68'
$ws.Range("E29").Value = '  +1.99%  '

Set-TextValue "D30" '0.999'
$ws.Range("E30").Value = '  -0.13%  '

Set-TextValue "D31" '7.26'
$ws.Range("E31").Value = '  +1.23%  '

Set-TextValue "D32" '2.08'
$ws.Range("E32").Value = '  -1.93%  '

Set-TextValue "D33" '27.67'
$ws.Range("E33").Value = '  +0.96%  '

$ws.Range("E34").Value = '  +2.21%  '

$ws.Range("D35").Value = '0.0₃0843'
$ws.Range("E35").Value = '  +5.79%  '

Set-TextValue "D36" '1.03'
$ws.Range("E36").Value = '  -0.78%  '

Set-TextValue "D37" '5.80'
$ws.Range("E37").Value = '  +0.76%  '

Set-TextValue "D38" '9.27'
$ws.Range("E38").Value = '  +3.40%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D39" '2.08'
$ws.Range("E39").Value = '  -2.49%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D40" '50.39'
$ws.Range("E40").Value = '  -0.05%  '

$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D41" '2.92'
$ws.Range("E41").Value = '  -0.75%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D42" '0.123'
$ws.Range("E42").Value = '  +9.54%  '

Set-TextValue "D43" '398.96'
$ws.Range("E43").Value = '  -5.55%  '

Set-TextValue "D44" '39.95'
$ws.Range("E44").Value = '  +4.81%  '

Set-TextValue "D45" '0.0355'
$ws.Range("E45").Value = '  +0.54%  '

Set-TextValue "D46" '0.271'
$ws.Range("E46").Value = '  -2.15%  '

$ws.Range("D47").Value = '2.719.92'
$ws.Range("E47").Value = '  -2.22%  '

Set-TextValue "D48" '133.20'
$ws.Range("E48").Value = '  +2.75%  '

Set-TextValue "D50" '0.108'
$ws.Range("E50").Value = '  -0.16%  '

Set-TextValue "D51" '2.17'
$ws.Range("E51").Value = '  +2.40%  '
